$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(6, 1).Value = 18
$ws.Cells.Item(11, 1).Value = 18
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(25, 1).Value = 18
$ws.Cells.Item(29, 1).Value = 0
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(41, 1).Value = 69
$ws.Cells.Item(44, 1).Value = 1
$ws.Cells.Item(47, 1).Value = 2
$ws.Cells.Item(48, 1).Value = 18
$ws.Cells.Item(53, 1).Value = 7
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(56, 1).Value = 0
$ws.Cells.Item(60, 1).Value = 9
$ws.Cells.Item(63, 1).Value = 0
$ws.Cells.Item(66, 1).Value = 2
$ws.Cells.Item(69, 1).Value = 0
$ws.Cells.Item(71, 1).Value = 4
$ws.Cells.Item(77, 1).Value = 9
$ws.Cells.Item(80, 1).Value = 9
$ws.Cells.Item(84, 1).Value = 36
$ws.Cells.Item(85, 1).Value = 4
$ws.Cells.Item(87, 1).Value = 9
$ws.Cells.Item(90, 1).Value = 4
$ws.Cells.Item(102, 1).Value = 2
$ws.Cells.Item(103, 1).Value = 2
$ws.Cells.Item(105, 1).Value = 18
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(117, 1).Value = 4
$ws.Cells.Item(121, 1).Value = 0
$ws.Cells.Item(123, 1).Value = 0
$ws.Cells.Item(125, 1).Value = 14
$ws.Cells.Item(131, 1).Value = 9
$ws.Cells.Item(133, 1).Value = 5
$ws.Cells.Item(134, 1).Value = 9
$ws.Cells.Item(138, 1).Value = 7
$ws.Cells.Item(145, 1).Value = 14
$ws.Cells.Item(153, 1).Value = 9
$ws.Cells.Item(160, 1).Value = 0
$ws.Cells.Item(164, 1).Value = 0
$ws.Cells.Item(165, 1).Value = 9
$ws.Cells.Item(172, 1).Value = 0
$ws.Cells.Item(173, 1).Value = 9
$ws.Cells.Item(174, 1).Value = 9
$ws.Cells.Item(179, 1).Value = 0
$ws.Cells.Item(184, 1).Value = 6
$ws.Cells.Item(187, 1).Value = 0
$ws.Cells.Item(189, 1).Value = 9
$ws.Cells.Item(192, 1).Value = 9
$ws.Cells.Item(200, 1).Value = 11
$ws.Cells.Item(201, 1).Value = 9
$ws.Cells.Item(203, 1).Value = 11
$ws.Cells.Item(205, 1).Value = 27
$ws.Cells.Item(209, 1).Value = 36
$ws.Cells.Item(211, 1).Value = 9
$ws.Cells.Item(213, 1).Value = 9
$ws.Cells.Item(214, 1).Value = 2
$ws.Cells.Item(215, 1).Value = 11
$ws.Cells.Item(217, 1).Value = 23
$ws.Cells.Item(218, 1).Value = 9
$ws.Cells.Item(219, 1).Value = 9
$ws.Cells.Item(221, 1).Value = 11
$ws.Cells.Item(223, 1).Value = 0
$ws.Cells.Item(225, 1).Value = 11
$ws.Cells.Item(226, 1).Value = 11
$ws.Cells.Item(231, 1).Value = 20
$ws.Cells.Item(232, 1).Value = 0
$ws.Cells.Item(233, 1).Value = 9
$ws.Cells.Item(234, 1).Value = 2
$ws.Cells.Item(236, 1).Value = 10
$ws.Cells.Item(237, 1).Value = 8
$ws.Cells.Item(244, 1).Value = 9
$ws.Cells.Item(247, 1).Value = 9
$ws.Cells.Item(250, 1).Value = 13
$ws.Cells.Item(251, 1).Value = 9
$ws.Cells.Item(256, 1).Value = 11
$ws.Cells.Item(262, 1).Value = 0
$ws.Cells.Item(264, 1).Value = 2
$ws.Cells.Item(275, 1).Value = 6
$ws.Cells.Item(276, 1).Value = 8
$ws.Cells.Item(277, 1).Value = 11
$ws.Cells.Item(278, 1).Value = 10
$ws.Cells.Item(280, 1).Value = 27
$ws.Cells.Item(283, 1).Value = 4
$ws.Cells.Item(286, 1).Value = 9
$ws.Cells.Item(287, 1).Value = 2
$ws.Cells.Item(288, 1).Value = 10
$ws.Cells.Item(295, 1).Value = 9
$ws.Cells.Item(302, 1).Value = 9
$ws.Cells.Item(303, 1).Value = 11
$ws.Cells.Item(306, 1).Value = 0
$ws.Cells.Item(307, 1).Value = 23
$ws.Cells.Item(308, 1).Value = 6
$ws.Cells.Item(309, 1).Value = 10
$ws.Cells.Item(312, 1).Value = 0
$ws.Cells.Item(313, 1).Value = 18
$ws.Cells.Item(314, 1).Value = 3
$ws.Cells.Item(316, 1).Value = 0
$ws.Cells.Item(317, 1).Value = 4
$ws.Cells.Item(319, 1).Value = 9
$ws.Cells.Item(321, 1).Value = 0
$ws.Cells.Item(324, 1).Value = 19
$ws.Cells.Item(325, 1).Value = 0
$ws.Cells.Item(327, 1).Value = 7
$ws.Cells.Item(330, 1).Value = 1
$ws.Cells.Item(331, 1).Value = 0
$ws.Cells.Item(334, 1).Value = 18
$ws.Cells.Item(337, 1).Value = 1
$ws.Cells.Item(339, 1).Value = 0
$ws.Cells.Item(340, 1).Value = 2
$ws.Cells.Item(350, 1).Value = 0
